$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45177 = 2023-09-08).
# Bump it by one day (-> 45178 = 2023-09-09) for every data row (2..490).
$lastRow = 490
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
